# Generate Report for Handoff
# Adds two new rows (a "45a89201-76de-4886-969a-bf523e71e3d7.md" file and a
# "b6c45c8e-3ec2-4425-8c5d-8dfa3008d792.md" file, both "Ready for handoff" /
# "Include") to the Overview sheet and to each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$uuid1 = "45a89201-76de-4886-969a-bf523e71e3d7"
$uuid2 = "b6c45c8e-3ec2-4425-8c5d-8dfa3008d792"

$hash1 = "82a893a4d8c3f93478dfebf9f362c59051800dd7"
$hash2 = "ca6a9567ed8f3b4c9c16ce804568f59394412795"

$commit1 = "9a0b43ddf51ae35a1d2a14df9b0b90d76a1d7e21"
$commit2 = "2c7e6ab0f1c9421dbb3a7cfbbb2240f6d5b0c4c3"
$hocommit1 = "4e2a5c1f3a6c48e5a9d0b6e2f7c3a8d1e5f4a9b2"
$hocommit2 = "7b1d4f2e9c6a3851d7f0c4b8a2e6d9f1c3a7b5e0"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "$uuid1.md"
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-26-11 22:26:52"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit1/e2e/$uuid1.md", "", "", "$uuid1.md")

$wsOverview.Range("A7").Value = "$uuid2.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-26-11 22:26:52"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit2/e2e/$uuid2.md", "", "", "$uuid2.md")

# ---------------------------------------------------------------------
# Locale sheets "zh-cn" / "de-de": same layout, different xlf + datetime.
# Columns: A Source File Name | B File Extension | C Status |
#          D Latest Handoff File | E Latest Handoff Datetime |
#          F Latest Target File | G Latest Handback File |
#          H Latest Handback DateTime | I Handoff Reason |
#          J Dependency From | K Error Detail
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Datetime6 = "2016-03-11 22:26:49"; Datetime7 = "2016-03-11 22:26:49" },
    @{ Name = "de-de"; Datetime6 = "2016-03-11 22:26:52"; Datetime7 = "2016-03-11 22:26:52" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)
    $loc = $locale.Name

    # --- row 6 : uuid1 ---
    $ws.Range("A6").Value = "$uuid1.md"
    $ws.Range("B6").Value = ".md"
    $ws.Range("C6").Value = "Ready for handoff"
    $ws.Range("D6").Value = "$uuid1.$hash1.$loc.xlf"
    $ws.Range("E6").Value = $locale.Datetime6
    $ws.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H6").Value = "0001-01-01 00:00:00"
    $ws.Range("I6").Value = "Include"

    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit1/e2e/$uuid1.md", "", "", "$uuid1.md")
    $ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hocommit1/ol-handoff/OpenLocalizationTestOrg/oltest.$loc/ci/ht/$uuid1.$hash1.$loc.xlf", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hocommit1/ol-handoff/OpenLocalizationTestOrg/oltest.$loc/ci/ht/$uuid1.$hash1.$loc.xlf", "", "", "$uuid1.$hash1.$loc.xlf")

    # --- row 7 : uuid2 ---
    $ws.Range("A7").Value = "$uuid2.md"
    $ws.Range("B7").Value = ".md"
    $ws.Range("C7").Value = "Ready for handoff"
    $ws.Range("D7").Value = "$uuid2.$hash2.$loc.xlf"
    $ws.Range("E7").Value = $locale.Datetime7
    $ws.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H7").Value = "0001-01-01 00:00:00"
    $ws.Range("I7").Value = "Include"

    $ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$commit2/e2e/$uuid2.md", "", "", "$uuid2.md")
    $ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hocommit2/ol-handoff/OpenLocalizationTestOrg/oltest.$loc/ci/ht/$uuid2.$hash2.$loc.xlf", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hocommit2/ol-handoff/OpenLocalizationTestOrg/oltest.$loc/ci/ht/$uuid2.$hash2.$loc.xlf", "", "", "$uuid2.$hash2.$loc.xlf")
}

Write-Host "Added handoff rows for $uuid1 and $uuid2 to Overview, zh-cn, de-de."
